$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clothingValues = @(
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Blouse",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Halter,Jumpsuit",
    "Caftan,Jumpsuit",
    "Trunks,Dress",
    "Jumpsuit,Sweatpants",
    "Jumpsuit,Parka",
    "Blouse,Jumpsuit",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Halter",
    "Halter,Jumpsuit",
    "Kaftan,Jumpsuit",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Blouse",
    "Jumpsuit,Blouse",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Blouse",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Blouse",
    "Jumpsuit,Trunks",
    "Trunks,Blouse",
    "Parka,Blouse",
    "Blouse,Sweatpants",
    "Blouse,Blazer",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Dress",
    "Jumpsuit,Dress",
    "Jumpsuit,Dress",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Blouse",
    "Halter,Jumpsuit",
    "Caftan,Trunks",
    "Dress,Jumpsuit",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Blouse",
    "Sweatpants,Parka",
    "Jumpsuit,Blouse",
    "Halter,Jumpsuit",
    "Blouse,Trunks",
    "Jumpsuit,Kaftan",
    "Sweatpants,Trunks",
    "Halter,Blouse",
    "Jumpsuit,Blouse",
    "Jumpsuit,Blouse",
    "Trunks,Blazer",
    "Jumpsuit,Dress",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Cutoffs",
    "Sweatpants,Trunks",
    "Jumpsuit,Dress",
    "Jumpsuit,Blouse",
    "Jumpsuit,Kaftan",
    "Caftan,Trunks",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Halter",
    "Caftan,Trunks",
    "Jumpsuit,Trunks",
    "Jumpsuit,Trunks",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Blouse",
    "Jumpsuit,Caftan",
    "Jumpsuit,Blouse",
    "Dress,Jumpsuit",
    "Dress,Jumpsuit",
    "Jumpsuit,Dress",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Dress",
    "Caftan,Parka",
    "Trunks,Jumpsuit",
    "Jumpsuit,Halter",
    "Caftan,Blouse",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Blouse",
    "Jumpsuit,Kaftan",
    "Halter,Jumpsuit",
    "Jumpsuit,Halter",
    "Jumpsuit,Kaftan",
    "Jumpsuit,Halter"
)

for ($i = 0; $i -lt $clothingValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $clothingValues[$i]
}
